# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 52 (shifting existing rows 52-63
# down to 53-64) on the "Arveja Verde" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 52; rows 52..63 shift down to 53..64.
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new record's data.
$ws.Range("A52").Value = 9
$ws.Range("B52").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 44463
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = 100112022
$ws.Range("G52").Value = "Arveja Verde"
$ws.Range("H52").Value = "Perfection"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 15
$ws.Range("K52").Value = 30000
$ws.Range("L52").Value = 31000
$ws.Range("M52").Value = 30533
$ws.Range("N52").Value = "`$/malla 25 kilos"
$ws.Range("O52").Value = "Provincia de Huasco"
$ws.Range("P52").Value = 1221
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
